$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet / sheet tab.
$ws.Name = "New as of 2022-08-24"

# 2. Insert a new row at position 6: this is a duplicate of the
#    "АБВГ.123456.009 / Крышка нижняя" row (row 5), pushing the old
#    rows 6 ("Крышка") and 7 ("Крышка верхняя") down by one.
$ws.Rows.Item(6).Insert()

$ws.Range("A6").Value = "АБВГ.123456.009"
$ws.Range("B6").Value = 7
$ws.Range("C6").Value = "Крышка нижняя"
$ws.Range("D6").Value = "Хоторн"
$ws.Range("E6").Value = "SpaceX"
$ws.Range("F6").Value = 44692

# 3. Append a new row 9: another "АБВГ.123456.011 / Крышка верхняя" entry
#    (duplicate of what is now row 8), with its own quantity and date.
$ws.Rows.Item(9).Insert()

$ws.Range("A9").Value = "АБВГ.123456.011"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = "Крышка верхняя"
$ws.Range("D9").Value = "Хоторн"
$ws.Range("E9").Value = "SpaceX"
$ws.Range("F9").Value = 44042
